$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "17.04.2023 12:45 (CET)"
$ws.Range("C7").Value = '{"ProposedVersion":"https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/4ed84f92a361098006df6448cc6932c63b60c9fc","UpdatedVersion":"https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/8c89276d5970a8c6cdbdf7ced361b6cdaa31cc08","Domain":"Process","LobeOwner":"member3","Result":"accept by lobe owner"}'
$ws.Range("D7").Value = "ef9f55d81e223f12e16df6d3336f9e28b03872bced4985157868c6f5f1b750be"
